# ---------------------------------------------------------------------------
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
#   * new "Player Info" sheet (player bio) inserted as the first tab
#   * existing sheet kept as "ODI Batting", MATCH_CARD_LINK column replaced
#     by a plain MATCH_CODE column (just the numeric code, not the full URL)
#   * new "ODI Batting Extra" sheet (per-match batting detail) appended
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- sheet bookkeeping (do ALL add/move/rename first, before grabbing the
#     worksheet references we will actually write data through - this COM
#     shim's worksheet handles track *position*, so resolving them by name
#     only after the final layout is in place keeps them stable) ----------
$wb.Worksheets.Item(1).Name = "ODI Batting"

$wsInfo = $wb.Worksheets.Add()
$wsInfo.Name = "Player Info"
$wsInfo.Move($wb.Worksheets.Item(1), $null)

$wsExtra = $wb.Worksheets.Add()
$wsExtra.Name = "ODI Batting Extra"
$wsExtra.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-resolve every sheet by name now that the final tab order is set.
$info = $wb.Worksheets.Item("Player Info")
$batting = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# =============================================================================
# "Player Info" sheet
# =============================================================================
$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle $info.Range("A1:D1")

# ID looks numeric ("4586") - force text so the value is kept as a string,
# matching how the source data is stored (inline string, not a number).
$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "4586"
$info.Range("B2").Value = "Lorcan John Tucker"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"

# =============================================================================
# "ODI Batting" sheet - MATCH_CARD_LINK -> MATCH_CODE
# =============================================================================
$batting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2="4284"; 3="4285"; 4="4343"; 5="4347"; 6="4352"; 7="4391"; 8="4394";
    9="4397"; 10="4426"; 11="4427"; 12="4428"; 13="4439"; 14="4442";
    15="4444"; 16="4446"; 17="4448"; 18="4466"; 19="4467"; 20="4468";
    21="4474"; 22="4475"; 23="4478"; 24="4492"; 25="4494"; 26="4496";
    27="4519"; 28="4605"; 29="4608"; 30="4614"; 31="4693"; 32="4694";
    33="4696"; 34="4726"; 35="4729"; 36="4734"
}

$batting.Range("D2:D36").NumberFormat = "@"
foreach ($row in $matchCodes.Keys) {
    $batting.Range("D$row").Value = $matchCodes[$row]
}

# Rows whose INNING_NUMBER (column B) cell was a blank placeholder get that
# cell removed entirely (no formatting, no content at all).
$blankInningRows = @(12, 19, 21, 22, 26, 33, 35)
foreach ($row in $blankInningRows) {
    $batting.Range("B$row").ClearContents()
}

# =============================================================================
# "ODI Batting Extra" sheet
# =============================================================================
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle $extra.Range("A1:F1")

# MATCH_CODE (A) and NUM_4/NUM_6/PERCENT (C/D/E) are numeric-looking text
# values in the source data - force text so they don't get silently
# coerced into numbers / percentages.
$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:E21").NumberFormat = "@"

$extraRows = @(
    @{ r=2;  code="4448"; pos=6;    n4="3"; n6="0"; pct="9.57%";  mom="NO" },
    @{ r=3;  code="4466"; pos=6;    n4="0"; n6="0"; pct="4.12%";  mom="NO" },
    @{ r=4;  code="4467"; pos=6;    n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=5;  code="4468"; pos=6;    n4="0"; n6="0"; pct="3.68%";  mom="NO" },
    @{ r=6;  code="4474"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=7;  code="4475"; pos=8;    n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=8;  code="4478"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=9;  code="4492"; pos=7;    n4="0"; n6="0"; pct="6.14%";  mom="NO" },
    @{ r=10; code="4494"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=11; code="4496"; pos=6;    n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=12; code="4519"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=13; code="4605"; pos=6;    n4="1"; n6="2"; pct="8.67%";  mom="NO" },
    @{ r=14; code="4608"; pos=6;    n4="0"; n6="1"; pct="8.80%";  mom="NO" },
    @{ r=15; code="4614"; pos=7;    n4="1"; n6="0"; pct="3.90%";  mom="NO" },
    @{ r=16; code="4693"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=17; code="4694"; pos=5;    n4="0"; n6="0"; pct="3.74%";  mom="NO" },
    @{ r=18; code="4696"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=19; code="4726"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=20; code="4729"; pos=$null; n4=$null; n6=$null; pct=$null; mom="NO" },
    @{ r=21; code="4734"; pos=5;    n4="4"; n6="0"; pct="27.72%"; mom="NO" }
)

foreach ($row in $extraRows) {
    $r = $row.r
    $extra.Range("A$r").Value = $row.code
    if ($row.pos -ne $null) {
        $extra.Range("B$r").Value = $row.pos
    }
    if ($row.n4 -ne $null) {
        $extra.Range("C$r").Value = $row.n4
        $extra.Range("D$r").Value = $row.n6
        $extra.Range("E$r").Value = $row.pct
    }
    $extra.Range("F$r").Value = $row.mom
}

Write-Output "done"
